$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.893.88"
$ws.Range("E2").Value = "'  -6.64%  "
$ws.Range("D3").Value = "'2.457.06"
$ws.Range("E3").Value = "'  -9.84%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "'  -0.25%  "
$ws.Range("D5").Value = "'469.18"
$ws.Range("E5").Value = "'  -7.15%  "
$ws.Range("D6").Value = "'133.40"
$ws.Range("E6").Value = "'  -5.57%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("D8").Value = "'0.497"
$ws.Range("E8").Value = "'  -6.51%  "
$ws.Range("D9").Value = "'2.452.63"
$ws.Range("E9").Value = "'  -10.41%  "
$ws.Range("D10").Value = "'0.0963"
$ws.Range("E10").Value = "'  -8.68%  "
$ws.Range("E11").Value = "'  -12.19%  "
$ws.Range("D12").Value = "'0.318"
$ws.Range("E12").Value = "'  -9.07%  "
$ws.Range("E13").Value = "'  -3.74%  "
$ws.Range("D14").Value = "'2.856.77"
$ws.Range("E14").Value = "'  -10.71%  "
$ws.Range("D15").Value = "'54.741.66"
$ws.Range("E15").Value = "'  -6.93%  "
$ws.Range("E16").Value = "'  -0.73%  "
$ws.Range("D17").Value = "'19.86"
$ws.Range("E17").Value = "'  -8.40%  "
$ws.Range("D18").Value = "'2.449.56"
$ws.Range("E18").Value = "'  -10.46%  "
$ws.Range("D19").Value = "'4.24"
$ws.Range("E19").Value = "'  -11.08%  "
$ws.Range("D20").Value = "'313.84"
$ws.Range("E20").Value = "'  -8.27%  "
$ws.Range("D21").Value = "'9.65"
$ws.Range("E21").Value = "'  -12.18%  "
$ws.Range("D22").Value = "'0.995"
$ws.Range("E22").Value = "'  -0.39%  "
$ws.Range("D23").Value = "'5.70"
$ws.Range("E23").Value = "'  +1.03%  "
$ws.Range("D24").Value = "'5.41"
$ws.Range("E24").Value = "'  -13.63%  "
$ws.Range("D25").Value = "'56.79"
$ws.Range("E25").Value = "'  -10.43%  "
$ws.Range("E26").Value = "'  +1.28%  "
$ws.Range("B27").Value = "'Polygon"
$ws.Range("C27").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").Value = "'0.389"
$ws.Range("E27").Value = "'  -9.02%  "
$ws.Range("B28").Value = "'Kaspa"
$ws.Range("C28").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.158"
$ws.Range("D29").Value = "'2.528.79"
$ws.Range("E29").Value = "'  -10.92%  "
$ws.Range("D30").Value = "'7.19"
$ws.Range("E30").Value = "'  -4.28%  "
$ws.Range("E31").Value = "'  -0.13%  "
$ws.Range("D32").Value = "'0.0₃0727"
$ws.Range("E32").Value = "'  -12.39%  "
$ws.Range("D33").Value = "'147.00"
$ws.Range("E33").Value = "'  -2.56%  "
$ws.Range("D34").Value = "'17.86"
$ws.Range("E34").Value = "'  -7.29%  "
$ws.Range("E35").Value = "'  -10.33%  "
$ws.Range("E36").Value = "'  -7.00%  "
$ws.Range("E37").Value = "'  -14.78%  "
$ws.Range("E38").Value = "'  -6.33%  "
$ws.Range("D39").Value = "'0.806"
$ws.Range("E39").Value = "'  -14.98%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "'  +0.24%  "
$ws.Range("E41").Value = "'  -8.74%  "
$ws.Range("D42").Value = "'0.599"
$ws.Range("E42").Value = "'  -0.53%  "
$ws.Range("E43").Value = "'  -6.26%  "
$ws.Range("E44").Value = "'  -8.84%  "
$ws.Range("E45").Value = "'  -10.35%  "
$ws.Range("D46").Value = "'10.09"
$ws.Range("E46").Value = "'  -2.64%  "
$ws.Range("D47").Value = "'1.943.23"
$ws.Range("E47").Value = "'  -11.41%  "
$ws.Range("E48").Value = "'  -0.73%  "
$ws.Range("D49").Value = "'0.0220"
$ws.Range("E49").Value = "'  -3.52%  "
$ws.Range("D50").Value = "'236.40"
$ws.Range("E50").Value = "'  +4.63%  "
$ws.Range("E51").Value = "'  -11.23%  "
